# Add data for 2022-06-02
# Updates the "through" date labels and refreshes the May / Total rows
# with the latest carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / update the "as of" date references
$ws.Name = "Through 2022-05-25"
$ws.Range("A6").Value = "May (through 05-25)"

# April (row 5) - only the 2022 column changes
$ws.Range("I5").Value = 116

# May (row 6) - updated counts for every year
$ws.Range("B6").Value = 15
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 48
$ws.Range("E6").Value = 40
$ws.Range("F6").Value = 35
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 94
$ws.Range("I6").Value = 93

# Total (row 7) - updated counts for every year
$ws.Range("B7").Value = 104
$ws.Range("C7").Value = 200
$ws.Range("D7").Value = 301
$ws.Range("E7").Value = 286
$ws.Range("F7").Value = 190
$ws.Range("G7").Value = 312
$ws.Range("H7").Value = 617
$ws.Range("I7").Value = 644
